# update new orleans xlsx files
#
# 1. hotel_info gains a new "State" column (inserted between Hotel_Name and
#    City) populated with "Louisiana" for the existing hotel row.
# 2. The worksheets are reordered so that review_info comes before
#    hotel_info.

$wb = $excel.ActiveWorkbook

$hotelInfo  = $wb.Worksheets.Item("hotel_info")
$reviewInfo = $wb.Worksheets.Item("review_info")

# Insert a new column before the existing "City" column (column C) and
# populate it with the hotel's state.
$hotelInfo.Columns.Item(3).Insert()
$hotelInfo.Range("C1").Value = "State"
$hotelInfo.Range("C2").Value = "Louisiana"

# Reorder the sheets: review_info first, hotel_info second.
$reviewInfo.Move($hotelInfo)
